$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column B, shifting existing B->D and C->E
$ws.Columns("B:C").Insert()

# Keep the (cosmetic) ~8-character column width consistent across C:E, matching
# the width that column C already had before the insert.
$ws.Range("C1:E1").ColumnWidth = 7.166666666666667

# New header row values for the freshly inserted columns
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# Fill the new B and C columns (rows 2-27) with the "UN" placeholder used elsewhere
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
}

# New data points for Jun_15 column
$ws.Range("C14").Value = '6/13/2018,Reiterates,Buy,$107.00'
$ws.Range("C27").Value = '6/14/2018,Reiterates,Buy,$108.00'
